$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Fix BT rotation bug: PheroLevel (C3) and MaxHp (J3) values for Onyscidus row
$ws.Range("C3").Value = 10
$ws.Range("J3").Value = 220

# Update the active selection to match the saved view state
$ws.Range("C4").Select()
